# Commit: "Inserçao da funcao de remocao de registro - Correcao de falhas
#          na coluna Data do Registro"
#
# 1) Add two new tracked measurements (columns J/K): "Percentual de Gordura
#    Corporal" and "Percentual de Musculatura Corporal", matching the header
#    style already used by the other columns.
# 2) Fix the "Data do Registro" column: row 15 (B15) was stamped with a
#    mismatched/buggy date style versus every other row in the column; give
#    it the same style as B2:B14. Also corrects the drifted timestamp in A15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header columns, carrying over the existing header formatting ---
$ws.Range("A1").Copy($ws.Range("J1"))
$ws.Range("J1").Value = "Percentual de Gordura Corporal"

$ws.Range("A1").Copy($ws.Range("K1"))
$ws.Range("K1").Value = "Percentual de Musculatura Corporal"

# --- Fix column B's inconsistent style on row 15 ---
$ws.Range("B2").Copy($ws.Range("B15"))
$ws.Range("B15").Value = 45770

# --- Correct the drifted "Data e Hora do Lancamento" timestamp on row 15 ---
$ws.Range("A15").Value = 45770.47045805556
